$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Item Actividad" column (C) for rows 11-15 with the tasks
# completed while finishing the gallery page
$ws.Range("C11").Value = "Mejora en diseño de navbar"
$ws.Range("C12").Value = "Mejora en diseño de footer"
$ws.Range("C13").Value = "Implementacion de imágenes y videos"
$ws.Range("C14").Value = "Recopilacion de informacion"
$ws.Range("C15").Value = "Ingreso de informacion relevante"

# Update the view so it shows the newly entered rows, selecting C16
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("C16").Select()
